$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Sending cluster ECs -> Target cluster MuSCs, plus refreshed TPM-derived metrics ---
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 3.895302
$ws.Range("H2").Value = 11.685906
$ws.Range("I2").Value = 0.1607797697193069
$ws.Range("J2").Value = 0.1607797697193069
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7146560000000001
$ws.Range("N2").Value = 2.143968
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2.783800946112
$ws.Range("R2").Value = 25.054208515008
$ws.Range("S2").Value = 0.1607797697193069
$ws.Range("T2").Value = 0.1607797697193069

# --- Row 3: Sending cluster ECs -> FAPs (target cluster MuSCs unchanged), refreshed metrics ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 9.844169000000001
$ws.Range("H3").Value = 29.532507
$ws.Range("I3").Value = 0.4063210567236994
$ws.Range("J3").Value = 0.4063210567236994
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7146560000000001
$ws.Range("N3").Value = 2.143968
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 7.035194440864001
$ws.Range("R3").Value = 63.316749967776
$ws.Range("S3").Value = 0.4063210567236994
$ws.Range("T3").Value = 0.4063210567236994

# --- Row 4: Sending cluster FAPs -> MuSCs, target cluster ECs -> MuSCs, refreshed metrics ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 10.488092
$ws.Range("H4").Value = 31.464276
$ws.Range("I4").Value = 0.4328991735569938
$ws.Range("J4").Value = 0.4328991735569938
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7146560000000001
$ws.Range("N4").Value = 2.143968
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 7.495377876352001
$ws.Range("R4").Value = 67.45840088716801
$ws.Range("S4").Value = 0.4328991735569938
$ws.Range("T4").Value = 0.4328991735569938

# --- Remove the now-duplicate rows 5-7 (collapsed into rows 2-4 above) ---
$ws.Rows("5:7").Delete()
